$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank rows 11-14 so the trailing rows (fix memory leaks, blank, Total)
# shift up from 15-17 to 11-13.
$ws.Range("A11:A14").EntireRow.Delete()

# Update point values
$ws.Range("B9").Value = 25
$ws.Range("B11").Value = 15

# Update the saved selection to match the new layout
$ws.Range("C12").Select()
